# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer update
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 08:58"

# Row 17 (e.g. country rank 21) - refreshed stats
$ws.Range("B17").Value = 287300
$ws.Range("C17").Value = 626
$ws.Range("D17").Value = 265215
$ws.Range("E17").Value = 15932
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 6153

# Row 56 (rank 60) - refreshed stats
$ws.Range("B56").Value = 41299
$ws.Range("C56").Value = 276
$ws.Range("D56").Value = 34164
$ws.Range("E56").Value = 6321
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 814

# Row 73 - refreshed stats
$ws.Range("D73").Value = 10290
$ws.Range("E73").Value = 11119

# Rows 146/147: Georgia overtakes Republica de Chipre in ranking, swapping position
$ws.Range("A146").Value = "Georgia"
$ws.Range("B146").Value = 1306
$ws.Range("C146").Value = 23
$ws.Range("D146").Value = 1085
$ws.Range("E146").Value = 204
$ws.Range("H146").Value = 17

$ws.Range("A147").Value = "Republica de Chipre"
$ws.Range("B147").Value = 1305
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 870
$ws.Range("E147").Value = 415
$ws.Range("H147").Value = 20

# Row 158 - refreshed stats
$ws.Range("D158").Value = 430
$ws.Range("E158").Value = 460

# Rows 213/214: Islas Malvinas and Montserrat swap position
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
